{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n// (plus the blank paragraph that separates it from the bibliography),\n// which the site's rebuild no longer emits.\nconst body = context.document.body;\n\nconst results = body.search(\"Ver no Jupiter Salvar em pdf Salvar em docx\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Target paragraph \"Ver no Jupiter ...\" not found.');\n}\n\n// Paragraph layout around the target:\n//   [empty paragraph]\n//   Ver no Jupiter Salvar em pdf Salvar em docx      <- targetPara\n//   \u00a9 2020 . Contact: ...                            <- nextPara\n// All three are removed; the paragraphs before/after the block stay.\nconst targetPara = results.items[0].paragraphs.getFirst();\nconst prevPara = targetPara.getPrevious();\nconst nextPara = targetPara.getNext();\n\n// Delete from the bottom up so earlier references stay valid.\nnextPara.delete();\ntargetPara.delete();\nprevPara.delete();\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Ver no Jupiter ...\" footer paragraph.\n$findRng = $d.Content\n$found = $findRng.Find.Execute(\"Ver no Jupiter Salvar em pdf Salvar em docx\")\nif (-not $found) {\n    throw \"Target paragraph 'Ver no Jupiter ...' not found.\"\n}\n\n# Rebuild a plain Range from the hit's Start/End so later Range operations\n# are not tied to the (stateful) Find range object.\n$hit = $d.Range($findRng.Start, $findRng.End)\n\n# Resolve which paragraph index contains the hit, scanning the live\n# collection instead of trusting a hard-coded number.\n$count = $d.Paragraphs.Count\n$targetIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $pr = $d.Paragraphs.Item($i).Range\n    if ($pr.Start -le $hit.Start -and $hit.End -le $pr.End) {\n        $targetIdx = $i\n        break\n    }\n}\nif ($targetIdx -eq -1) {\n    throw \"Could not resolve paragraph index for the match.\"\n}\n\n# Layout around the target paragraph:\n#   targetIdx - 1 : blank separator paragraph\n#   targetIdx     : \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   targetIdx + 1 : \"\u00a9 2020 . Contact: ...\" copyright/footer line\n# Remove all three, keeping the paragraphs before/after the block intact.\n$startRng = $d.Paragraphs.Item($targetIdx - 1).Range\n$endRng = $d.Paragraphs.Item($targetIdx + 1).Range\n$block = $d.Range($startRng.Start, $endRng.End)\n$block.Delete()\n"}
